$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.859238885935099
$ws.Range("C2").Value = 0.3418508317708984
$ws.Range("D2").Value = 0.03885608809627428
$ws.Range("E2").Value = 0.087184406597693
$ws.Range("F2").Value = 1.557330138999873
$ws.Range("I2").Value = 0.9864072933335706
$ws.Range("M2").Value = 0.3925723809337143
$ws.Range("B3").Value = 0.7684864284123023
$ws.Range("C3").Value = 0.3006474811030273
$ws.Range("D3").Value = 0.03843379391274127
$ws.Range("E3").Value = 0.08181199126968863
$ws.Range("F3").Value = 1.488194764439896
$ws.Range("I3").Value = 0.956434105413777
$ws.Range("M3").Value = 0.3545240259354614
$ws.Range("B4").Value = 0.7131505995822636
$ws.Range("C4").Value = 0.275457126316013
$ws.Range("D4").Value = 0.03819084857061128
$ws.Range("E4").Value = 0.07857295233063155
$ws.Range("F4").Value = 1.446716133549685
$ws.Range("I4").Value = 0.938591943007097
$ws.Range("M4").Value = 0.331377648331717
$ws.Range("B5").Value = 0.69069678383309
$ws.Range("C5").Value = 0.2652181449358864
$ws.Range("D5").Value = 0.03809587335061693
$ws.Range("E5").Value = 0.0772677193690221
$ws.Range("F5").Value = 1.430054052876926
$ws.Range("I5").Value = 0.931460468732908
$ws.Range("M5").Value = 0.3219984755833636
$ws.Range("B6").Value = 0.6869741073581679
$ws.Range("C6").Value = 0.263519528058282
$ws.Range("D6").Value = 0.03808034358176471
$ws.Range("E6").Value = 0.07705186575601886
$ws.Range("F6").Value = 1.427301774360572
$ws.Range("I6").Value = 0.9302846556908548
$ws.Range("M6").Value = 0.3204442527094145
$ws.Range("B7").Value = 0.7128473925129413
$ws.Range("C7").Value = 0.2753189347185128
$ws.Range("D7").Value = 0.03818955150680026
$ws.Range("E7").Value = 0.07855529036418218
$ws.Range("F7").Value = 1.446490452300424
$ws.Range("I7").Value = 0.9384952036657239
$ws.Range("M7").Value = 0.3312509435946609
$ws.Range("B8").Value = 0.8278664377278346
$ws.Range("C8").Value = 0.3276207234189883
$ws.Range("D8").Value = 0.03870703842833834
$ws.Range("E8").Value = 0.08531943730448432
$ws.Range("F8").Value = 1.533288758455683
$ws.Range("I8").Value = 0.9759548793168022
$ws.Range("M8").Value = 0.3794079595848103
$ws.Range("B9").Value = 1.056554149964825
$ws.Range("C9").Value = 0.4310988719794864
$ws.Range("D9").Value = 0.03985532292252003
$ws.Range("E9").Value = 0.09907132646738148
$ws.Range("F9").Value = 1.711361868640665
$ws.Range("I9").Value = 1.05395580332636
$ws.Range("M9").Value = 0.4756047487397979
$ws.Range("B10").Value = 1.22659756559608
$ws.Range("C10").Value = 0.5077623236036857
$ws.Range("D10").Value = 0.04078569740477889
$ws.Range("E10").Value = 0.1094927007643633
$ws.Range("F10").Value = 1.847220030664971
$ws.Range("I10").Value = 1.114159691422884
$ws.Range("M10").Value = 0.5474332886716837
$ws.Range("B11").Value = 1.304421482805481
$ws.Range("C11").Value = 0.542795228998898
$ws.Range("D11").Value = 0.041228946823054
$ws.Range("E11").Value = 0.1143072629544832
$ws.Range("F11").Value = 1.910169149029684
$ws.Range("I11").Value = 1.142204973610276
$ws.Range("M11").Value = 0.580378131204796
$ws.Range("B12").Value = 1.333960866045459
$ws.Range("C12").Value = 0.5560853916534256
$ws.Range("D12").Value = 0.0413997665543917
$ws.Range("E12").Value = 0.1161413861725151
$ws.Range("F12").Value = 1.934175194312047
$ws.Range("I12").Value = 1.152921812438521
$ws.Range("M12").Value = 0.5928935943051954
$ws.Range("B13").Value = 1.327595930807661
$ws.Range("C13").Value = 0.5532220323207184
$ws.Range("D13").Value = 0.04136284386335376
$ws.Range("E13").Value = 0.1157458823683299
$ws.Range("F13").Value = 1.928997505798861
$ws.Range("I13").Value = 1.150609415933047
$ws.Range("M13").Value = 0.5901963687232978
$ws.Range("B14").Value = 1.306850310752168
$ws.Range("C14").Value = 0.5438881315094477
$ws.Range("D14").Value = 0.04124294019381125
$ws.Range("E14").Value = 0.1144579359139897
$ws.Range("F14").Value = 1.912140740839021
$ws.Range("I14").Value = 1.143084704606537
$ws.Range("M14").Value = 0.5814069780152806
$ws.Range("B15").Value = 1.294152071614974
$ws.Range("C15").Value = 0.5381740042016077
$ws.Range("D15").Value = 0.04116988537470689
$ws.Range("E15").Value = 0.1136704680072071
$ws.Range("F15").Value = 1.901837550442167
$ws.Range("I15").Value = 1.138488255788474
$ws.Range("M15").Value = 0.5760284628548931
$ws.Range("B16").Value = 1.22152118256065
$ws.Range("C16").Value = 0.505476128241412
$ws.Range("D16").Value = 0.04075714068718383
$ws.Range("E16").Value = 0.1091795720661253
$ws.Range("F16").Value = 1.843129545314468
$ws.Range("I16").Value = 1.11234028605476
$ws.Range("M16").Value = 0.5452857932367721
$ws.Range("B17").Value = 1.177086010770722
$ws.Range("C17").Value = 0.485458432132873
$ws.Range("D17").Value = 0.04050912755174352
$ws.Range("E17").Value = 0.1064436997860199
$ws.Range("F17").Value = 1.807410370923634
$ws.Range("I17").Value = 1.096469392312201
$ws.Range("M17").Value = 0.5264960689365239
$ws.Range("B18").Value = 1.151572248965238
$ws.Range("C18").Value = 0.4739596112440836
$ws.Range("D18").Value = 0.04036835552143714
$ws.Range("E18").Value = 0.1048770269031962
$ws.Range("F18").Value = 1.786973341977557
$ws.Range("I18").Value = 1.087402706393007
$ws.Range("M18").Value = 0.5157140274994703
$ws.Range("B19").Value = 1.142941268937193
$ws.Range("C19").Value = 0.4700688188494269
$ws.Range("D19").Value = 0.04032101237579866
$ws.Range("E19").Value = 0.1043477581770205
$ws.Range("F19").Value = 1.780072099668502
$ws.Range("I19").Value = 1.084343440057268
$ws.Range("M19").Value = 0.5120677249783654
$ws.Range("B20").Value = 1.181811627981801
$ws.Range("C20").Value = 0.487587806782301
$ws.Range("D20").Value = 0.04053533388341179
$ws.Range("E20").Value = 0.1067342185639006
$ws.Range("F20").Value = 1.811201565235365
$ws.Range("I20").Value = 1.098152461867784
$ws.Range("M20").Value = 0.5284936390555401
$ws.Range("B21").Value = 1.312941912985366
$ws.Range("C21").Value = 0.5466290644647529
$ws.Range("D21").Value = 0.04127807744171719
$ws.Range("E21").Value = 0.1148359370380376
$ws.Range("F21").Value = 1.917087376346956
$ws.Range("I21").Value = 1.145292253933093
$ws.Range("M21").Value = 0.5839875408696145
$ws.Range("B22").Value = 1.399046968058201
$ws.Range("C22").Value = 0.5853561797363795
$ws.Range("D22").Value = 0.04178086281171289
$ws.Range("E22").Value = 0.1201948807503186
$ws.Range("F22").Value = 1.987274241148157
$ws.Range("I22").Value = 1.176665433310617
$ws.Range("M22").Value = 0.6204895364001715
$ws.Range("B23").Value = 1.353053611058556
$ws.Range("C23").Value = 0.5646735557015745
$ws.Range("D23").Value = 0.04151089767808003
$ws.Range("E23").Value = 0.1173287428475973
$ws.Range("F23").Value = 1.949722884996618
$ws.Range("I23").Value = 1.159868634851719
$ws.Range("M23").Value = 0.6009859683627212
$ws.Range("B24").Value = 1.179675074980537
$ws.Range("C24").Value = 0.4866250866314203
$ws.Range("D24").Value = 0.04052348035692432
$ws.Range("E24").Value = 0.1066028556867593
$ws.Range("F24").Value = 1.809487260072075
$ws.Range("I24").Value = 1.097391366546461
$ws.Range("M24").Value = 0.5275904740860966
$ws.Range("B25").Value = 0.9943406728968966
$ws.Range("C25").Value = 0.4029991464112754
$ws.Range("D25").Value = 0.03952983979711888
$ws.Range("E25").Value = 0.09529676884175586
$ws.Range("F25").Value = 1.662322545892437
$ws.Range("I25").Value = 1.032355068990114
$ws.Range("M25").Value = 0.4493843379078157

Write-Host "Applied 380 kV case update"
